$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 183, shifting existing rows 183-210 down to 184-211
$ws.Rows("183:183").Insert()

# Populate the newly inserted row 183 with its data
$ws.Range("A183").Value = 6
$ws.Range("B183").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C183").Value = "Metropolitana"
$ws.Range("D183").Value = 44474
$ws.Range("E183").Value = 13
$ws.Range("F183").Value = 100112032
$ws.Range("G183").Value = "Zapallo italiano"
$ws.Range("H183").Value = "Sin especificar"
$ws.Range("I183").Value = "Primera"
$ws.Range("J183").Value = 130
$ws.Range("K183").Value = 22000
$ws.Range("L183").Value = 24000
$ws.Range("M183").Value = 22769
$ws.Range("N183").Value = "$/caja 60 unidades"
$ws.Range("O183").Value = "Limache"
$ws.Range("P183").Value = 379
$ws.Range("Q183").Value = 60
$ws.Range("R183").Value = "Hortaliza"
